$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Redact/rename the query headers: "Consulta 1" -> "Consulta A", "Consulta 2" -> "Consulta B"
$ws.Range("A1").Value = "Consulta A"
$ws.Range("A9").Value = "Consulta B"

# Move the active selection to match the final saved view state
$ws.Range("G7").Select()
